$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "file name" column (A) previously held full Windows paths to the interact-*.pep.xml
# PhilosopherPipeline outputs; replace each block of rows (one block per raw file / 16 TMT
# channels) with the corresponding bare *.raw acquisition file name.
$fileNameBlocks = @(
    @{ Start = 2; End = 17; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H14_100pg_AGC300_1.raw" }
    @{ Start = 19; End = 34; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H14_100pg_AGC300_2.raw" }
    @{ Start = 36; End = 51; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H42_100pg_AGC300_1.raw" }
    @{ Start = 53; End = 68; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H42_100pg_AGC300_2.raw" }
    @{ Start = 70; End = 85; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H98_100pg_AGC300_1.raw" }
    @{ Start = 87; End = 102; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H98_100pg_AGC300_2.raw" }
    @{ Start = 104; End = 119; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H210_100pg_AGC300_1.raw" }
    @{ Start = 121; End = 136; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H210_100pg_AGC300_2.raw" }
    @{ Start = 138; End = 153; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H434_100pg_AGC300_1.raw" }
    @{ Start = 155; End = 170; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H434_100pg_AGC300_2.raw" }
    @{ Start = 172; End = 187; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_No126_100pg_AGC300_1.raw" }
    @{ Start = 189; End = 204; Value = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_No126_100pg_AGC300_2.raw" }
)

foreach ($block in $fileNameBlocks) {
    for ($r = $block.Start; $r -le $block.End; $r++) {
        $ws.Range("A$r").Value = $block.Value
    }
}

# Column A is much narrower now that it holds short file names instead of long paths;
# resize it to fit the new content (closest achievable grid value to bestFit ~63.78).
$ws.Columns.Item(1).ColumnWidth = 63

# Restore the active selection recorded in the workbook to cell A6.
$ws.Range("A6").Select() | Out-Null
